$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ebi3"
$ws.Range("C2").Value = "Il27ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 2.878032333333334
$ws.Range("H2").Value = 8.634097
$ws.Range("I2").Value = 0.3279446910817746
$ws.Range("J2").Value = 0.3279446910817746
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.087864
$ws.Range("N2").Value = 3.263592
$ws.Range("O2").Value = 0.3655718228968423
$ws.Range("P2").Value = 0.3655718228968423
$ws.Range("Q2").Value = 3.130907766269333
$ws.Range("R2").Value = 28.178169896424
$ws.Range("S2").Value = 0.1198873385281062
$ws.Range("T2").Value = 0.1198873385281062

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ebi3"
$ws.Range("C3").Value = "Il27ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 2.878032333333334
$ws.Range("H3").Value = 8.634097
$ws.Range("I3").Value = 0.3279446910817746
$ws.Range("J3").Value = 0.3279446910817746
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.742815333333333
$ws.Range("N3").Value = 5.228446
$ws.Range("O3").Value = 0.5856652838766928
$ws.Range("P3").Value = 0.5856652838766928
$ws.Range("Q3").Value = 5.015878880362445
$ws.Range("R3").Value = 45.142909923262
$ws.Range("S3").Value = 0.1920658205982619
$ws.Range("T3").Value = 0.1920658205982619

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ebi3"
$ws.Range("C4").Value = "Il27ra"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 2.878032333333334
$ws.Range("H4").Value = 8.634097
$ws.Range("I4").Value = 0.3279446910817746
$ws.Range("J4").Value = 0.3279446910817746
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.145108
$ws.Range("N4").Value = 0.435324
$ws.Range("O4").Value = 0.04876289322646489
$ws.Range("P4").Value = 0.04876289322646488
$ws.Range("Q4").Value = 0.4176255158253334
$ws.Range("R4").Value = 3.758629642428001
$ws.Range("S4").Value = 0.01599153195540659
$ws.Range("T4").Value = 0.01599153195540658

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Ebi3"
$ws.Range("C5").Value = "Il27ra"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.897936333333334
$ws.Range("H5").Value = 17.693809
$ws.Range("I5").Value = 0.6720553089182254
$ws.Range("J5").Value = 0.6720553089182254
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.087864
$ws.Range("N5").Value = 3.263592
$ws.Range("O5").Value = 0.3655718228968423
$ws.Range("P5").Value = 0.3655718228968423
$ws.Range("Q5").Value = 6.416152611325334
$ws.Range("R5").Value = 57.74537350192801
$ws.Range("S5").Value = 0.2456844843687361
$ws.Range("T5").Value = 0.2456844843687361

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Ebi3"
$ws.Range("C6").Value = "Il27ra"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.897936333333334
$ws.Range("H6").Value = 17.693809
$ws.Range("I6").Value = 0.6720553089182254
$ws.Range("J6").Value = 0.6720553089182254
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.742815333333333
$ws.Range("N6").Value = 5.228446
$ws.Range("O6").Value = 0.5856652838766928
$ws.Range("P6").Value = 0.5856652838766928
$ws.Range("Q6").Value = 10.27901387675711
$ws.Range("R6").Value = 92.511124890814
$ws.Range("S6").Value = 0.393599463278431
$ws.Range("T6").Value = 0.393599463278431

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ebi3"
$ws.Range("C7").Value = "Il27ra"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.897936333333334
$ws.Range("H7").Value = 17.693809
$ws.Range("I7").Value = 0.6720553089182254
$ws.Range("J7").Value = 0.6720553089182254
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.145108
$ws.Range("N7").Value = 0.435324
$ws.Range("O7").Value = 0.04876289322646489
$ws.Range("P7").Value = 0.04876289322646488
$ws.Range("Q7").Value = 0.8558377454573335
$ws.Range("R7").Value = 7.702539709116001
$ws.Range("S7").Value = 0.0327713612710583
$ws.Range("T7").Value = 0.0327713612710583
